# Scheduled market-data refresh: overwrite the computed price/profit
# columns (H:N) for the affected Leve rows on each job sheet with the
# latest pulled values. Mirrors chore: update Sheets via scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2425.5
$ws.Range("I41").Value = 2217
$ws.Range("J41").Value = 2634
$ws.Range("K41").Value = 2217
$ws.Range("L41").Value = 2634
$ws.Range("M41").Value = -1777
$ws.Range("N41").Value = -3514

$ws.Range("H106").Value = 5514.2
$ws.Range("I106").Value = 5514.2
$ws.Range("K106").Value = 5514.2
$ws.Range("M106").Value = -4883.2

$ws.Range("H112").Value = 1864.3914
$ws.Range("I112").Value = 1400
$ws.Range("J112").Value = 1885.5
$ws.Range("K112").Value = 4200
$ws.Range("L112").Value = 5656.5
$ws.Range("M112").Value = -3092
$ws.Range("N112").Value = -7872.5

$ws.Range("H116").Value = 27335.834
$ws.Range("I116").Value = 40615.266
$ws.Range("J116").Value = 5203.4443
$ws.Range("K116").Value = 40615.266
$ws.Range("L116").Value = 5203.4443
$ws.Range("M116").Value = -37173.266
$ws.Range("N116").Value = -12087.4443

$ws.Range("H129").Value = 1714.8235
$ws.Range("I129").Value = 1476.8
$ws.Range("K129").Value = 4430.4
$ws.Range("M129").Value = 569.6000000000004

$ws.Range("H132").Value = 149040.56
$ws.Range("J132").Value = 4587.8
$ws.Range("L132").Value = 13763.4
$ws.Range("N132").Value = -18823.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4162.467
$ws.Range("J74").Value = 4712.3335
$ws.Range("L74").Value = 4712.3335
$ws.Range("N74").Value = -6460.3335

$ws.Range("H77").Value = 4162.467
$ws.Range("J77").Value = 4712.3335
$ws.Range("L77").Value = 23561.6675
$ws.Range("N77").Value = -32297.6675

$ws.Range("H97").Value = 1804.5714
$ws.Range("J97").Value = 1168.75
$ws.Range("L97").Value = 1168.75
$ws.Range("N97").Value = -2160.75

$ws.Range("H102").Value = 10463.556
$ws.Range("I102").Value = 10463.556
$ws.Range("K102").Value = 10463.556
$ws.Range("M102").Value = -8841.556

$ws.Range("H124").Value = 34988.5
$ws.Range("J124").Value = 34988.5
$ws.Range("L124").Value = 34988.5
$ws.Range("N124").Value = -44808.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1200.8928
$ws.Range("I86").Value = 1115.6111
$ws.Range("K86").Value = 1115.6111
$ws.Range("M86").Value = 7.388899999999921

$ws.Range("H89").Value = 1200.8928
$ws.Range("I89").Value = 1115.6111
$ws.Range("K89").Value = 5578.0555
$ws.Range("M89").Value = 37.94449999999961

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H22").Value = 856.05
$ws.Range("J22").Value = 1945.3334
$ws.Range("L22").Value = 1945.3334
$ws.Range("N22").Value = -2645.3334

$ws.Range("H25").Value = 12837.5
$ws.Range("J25").Value = 50000
$ws.Range("L25").Value = 50000
$ws.Range("N25").Value = -50348

$ws.Range("H31").Value = 6562.4883
$ws.Range("I31").Value = 5832.304
$ws.Range("J31").Value = 7402.2
$ws.Range("K31").Value = 5832.304
$ws.Range("L31").Value = 7402.2
$ws.Range("M31").Value = -5537.304
$ws.Range("N31").Value = -7992.2

$ws.Range("H34").Value = 6562.4883
$ws.Range("I34").Value = 5832.304
$ws.Range("J34").Value = 7402.2
$ws.Range("K34").Value = 5832.304
$ws.Range("L34").Value = 7402.2
$ws.Range("M34").Value = -5630.304
$ws.Range("N34").Value = -7806.2

$ws.Range("H107").Value = 716.44446
$ws.Range("I107").Value = 653.94116
$ws.Range("J107").Value = 822.7
$ws.Range("K107").Value = 653.94116
$ws.Range("L107").Value = 822.7
$ws.Range("M107").Value = 1266.05884
$ws.Range("N107").Value = -4662.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 5255.75
$ws.Range("I115").Value = 614
$ws.Range("J115").Value = 9897.5
$ws.Range("K115").Value = 1842
$ws.Range("L115").Value = 29692.5
$ws.Range("M115").Value = -667
$ws.Range("N115").Value = -32042.5

$ws.Range("H131").Value = 933.5714
$ws.Range("I131").Value = 551.3226
$ws.Range("K131").Value = 1653.9678
$ws.Range("M131").Value = 3386.0322

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1429.7858
$ws.Range("I97").Value = 1486.3
$ws.Range("J97").Value = 1288.5
$ws.Range("K97").Value = 1486.3
$ws.Range("L97").Value = 1288.5
$ws.Range("M97").Value = -990.3
$ws.Range("N97").Value = -2280.5

$ws.Range("H132").Value = 3244.697
$ws.Range("I132").Value = 2040.1111
$ws.Range("K132").Value = 6120.3333
$ws.Range("M132").Value = -3590.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1243.5
$ws.Range("J22").Value = 1552.6364
$ws.Range("L22").Value = 1552.6364
$ws.Range("N22").Value = -2142.6364

$ws.Range("H27").Value = 1243.5
$ws.Range("J27").Value = 1552.6364
$ws.Range("L27").Value = 1552.6364
$ws.Range("N27").Value = -1766.6364

$ws.Range("H46").Value = 1914.5834
$ws.Range("J46").Value = 1886.2222
$ws.Range("L46").Value = 1886.2222
$ws.Range("N46").Value = -2262.2222

$ws.Range("H55").Value = 305.5
$ws.Range("I55").Value = 317.2857
$ws.Range("J55").Value = 295.1875
$ws.Range("K55").Value = 317.2857
$ws.Range("L55").Value = 295.1875
$ws.Range("M55").Value = -144.2857
$ws.Range("N55").Value = -641.1875

$ws.Range("H93").Value = 4300
$ws.Range("I93").Value = 3950
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 3950
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -2702
$ws.Range("N93").Value = -7496

$ws.Range("H100").Value = 100004
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 100004
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 100004
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -101086

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2890
$ws.Range("I107").Value = 1187.8889
$ws.Range("K107").Value = 3563.6667
$ws.Range("M107").Value = -1643.6667

$ws.Range("H132").Value = 2445.5103
$ws.Range("I132").Value = 2317.0698
$ws.Range("K132").Value = 6951.209400000001
$ws.Range("M132").Value = -4421.209400000001

$ws.Range("H136").Value = 2925.2666
$ws.Range("I136").Value = 2634.1304
$ws.Range("J136").Value = 3881.8572
$ws.Range("K136").Value = 7902.3912
$ws.Range("L136").Value = 11645.5716
$ws.Range("M136").Value = -5352.3912
